# Apply numeric value corrections to the Leve profit tables across all job sheets.
# Values come from an updated Market Board price snapshot; only raw data cells change.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 316679.44
$ws.Range("J17").Value = 332463.4
$ws.Range("L17").Value = 997390.2000000001
$ws.Range("N17").Value = -997726.2000000001
$ws.Range("H19").Value = 999
$ws.Range("J19").Value = 999
$ws.Range("L19").Value = 999
$ws.Range("N19").Value = -1349
$ws.Range("H33").Value = 2660538.8
$ws.Range("I33").Value = 3856411.8
$ws.Range("J33").Value = 3043
$ws.Range("K33").Value = 3856411.8
$ws.Range("L33").Value = 3043
$ws.Range("M33").Value = -3856182.8
$ws.Range("N33").Value = -3501
$ws.Range("H40").Value = 2039.8
$ws.Range("J40").Value = 2299.8572
$ws.Range("L40").Value = 2299.8572
$ws.Range("N40").Value = -2649.8572
$ws.Range("H97").Value = 6141.2
$ws.Range("I97").Value = 200
$ws.Range("J97").Value = 7626.5
$ws.Range("K97").Value = 600
$ws.Range("L97").Value = 22879.5
$ws.Range("M97").Value = -104
$ws.Range("N97").Value = -23871.5
$ws.Range("H113").Value = 4143.5713
$ws.Range("I113").Value = 4143.5713
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4143.5713
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -889.5712999999996
$ws.Range("H138").Value = 4529.937
$ws.Range("I138").Value = 5213.778
$ws.Range("J138").Value = 4112.6777
$ws.Range("K138").Value = 15641.334
$ws.Range("L138").Value = 12338.0331
$ws.Range("M138").Value = -10501.334
$ws.Range("N138").Value = -22618.0331
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 1727.5
$ws.Range("I30").Value = 1750
$ws.Range("J30").Value = 1705
$ws.Range("K30").Value = 1750
$ws.Range("L30").Value = 1705
$ws.Range("M30").Value = -1600
$ws.Range("N30").Value = -2005
$ws.Range("H33").Value = 6999.3335
$ws.Range("I33").Value = 6999.3335
$ws.Range("K33").Value = 6999.3335
$ws.Range("M33").Value = -6670.3335
$ws.Range("H40").Value = 99999
$ws.Range("J40").Value = 99999
$ws.Range("L40").Value = 99999
$ws.Range("N40").Value = -100297
$ws.Range("H45").Value = 55736.105
$ws.Range("I45").Value = 85101.586
$ws.Range("K45").Value = 85101.586
$ws.Range("M45").Value = -84724.586
$ws.Range("H61").Value = 1551456.5
$ws.Range("I61").Value = 3738.6453
$ws.Range("K61").Value = 3738.6453
$ws.Range("M61").Value = -3526.6453
$ws.Range("H74").Value = 1063671
$ws.Range("I74").Value = 2138
$ws.Range("K74").Value = 2138
$ws.Range("M74").Value = -1264
$ws.Range("H77").Value = 1063671
$ws.Range("I77").Value = 2138
$ws.Range("K77").Value = 10690
$ws.Range("M77").Value = -6322
$ws.Range("H136").Value = 1551456.5
$ws.Range("I136").Value = 3738.6453
$ws.Range("K136").Value = 11215.9359
$ws.Range("M136").Value = -8665.9359

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 844.625
$ws.Range("I22").Value = 830.3333
$ws.Range("K22").Value = 830.3333
$ws.Range("M22").Value = -657.3333
$ws.Range("H105").Value = 9089.111000000001
$ws.Range("I105").Value = 9123.462
$ws.Range("J105").Value = 8999.799999999999
$ws.Range("K105").Value = 9123.462
$ws.Range("L105").Value = 8999.799999999999
$ws.Range("M105").Value = -7376.462
$ws.Range("N105").Value = -12493.8
$ws.Range("H115").Value = 30000
$ws.Range("J115").Value = 30000
$ws.Range("L115").Value = 30000
$ws.Range("N115").Value = -33134

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1015.7778
$ws.Range("J22").Value = 799
$ws.Range("L22").Value = 799
$ws.Range("N22").Value = -1499
$ws.Range("H29").Value = 950
$ws.Range("J29").Value = 950
$ws.Range("L29").Value = 950
$ws.Range("N29").Value = -1530
$ws.Range("H31").Value = 2721.3455
$ws.Range("I31").Value = 2498.5186
$ws.Range("K31").Value = 2498.5186
$ws.Range("M31").Value = -2203.5186
$ws.Range("H33").Value = 5657.8
$ws.Range("I33").Value = 1196.3334
$ws.Range("J33").Value = 12350
$ws.Range("K33").Value = 1196.3334
$ws.Range("L33").Value = 12350
$ws.Range("M33").Value = -817.3334
$ws.Range("N33").Value = -13108
$ws.Range("H34").Value = 2721.3455
$ws.Range("I34").Value = 2498.5186
$ws.Range("K34").Value = 2498.5186
$ws.Range("M34").Value = -2296.5186
$ws.Range("H41").Value = 10000
$ws.Range("I41").Value = 10000
$ws.Range("K41").Value = 10000
$ws.Range("M41").Value = -9572
$ws.Range("H60").Value = 45000
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("H86").Value = 14121.692
$ws.Range("I86").Value = 9898
$ws.Range("J86").Value = 15998.889
$ws.Range("K86").Value = 9898
$ws.Range("L86").Value = 15998.889
$ws.Range("M86").Value = -8775
$ws.Range("N86").Value = -18244.889
$ws.Range("H89").Value = 14121.692
$ws.Range("I89").Value = 9898
$ws.Range("J89").Value = 15998.889
$ws.Range("K89").Value = 49490
$ws.Range("L89").Value = 79994.44499999999
$ws.Range("M89").Value = -43874
$ws.Range("N89").Value = -91226.44499999999
$ws.Range("H122").Value = 3173.5
$ws.Range("I122").Value = 3782.2222
$ws.Range("K122").Value = 11346.6666
$ws.Range("M122").Value = -8896.6666
$ws.Range("M60").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 78.2
$ws.Range("I2").Value = 93.73333
$ws.Range("J2").Value = 31.6
$ws.Range("K2").Value = 562.3999799999999
$ws.Range("L2").Value = 189.6
$ws.Range("M2").Value = -449.3999799999999
$ws.Range("N2").Value = -415.6
$ws.Range("H7").Value = 252.5
$ws.Range("I7").Value = 262.85715
$ws.Range("J7").Value = 180
$ws.Range("K7").Value = 788.5714499999999
$ws.Range("L7").Value = 540
$ws.Range("M7").Value = -676.5714499999999
$ws.Range("N7").Value = -764
$ws.Range("H87").Value = 16986.309
$ws.Range("I87").Value = 4199.4
$ws.Range("K87").Value = 12598.2
$ws.Range("M87").Value = -11350.2
$ws.Range("H90").Value = 16986.309
$ws.Range("I90").Value = 4199.4
$ws.Range("K90").Value = 37794.6
$ws.Range("M90").Value = -31554.6
$ws.Range("H121").Value = 1276.0769
$ws.Range("I121").Value = 399.1111
$ws.Range("J121").Value = 3249.25
$ws.Range("K121").Value = 1197.3333
$ws.Range("L121").Value = 9747.75
$ws.Range("M121").Value = 112.6667
$ws.Range("N121").Value = -12367.75
$ws.Range("H122").Value = 2756151
$ws.Range("I122").Value = 6061425.5
$ws.Range("K122").Value = 54552829.5
$ws.Range("M122").Value = -54550379.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 115.64286
$ws.Range("I2").Value = 104.2
$ws.Range("J2").Value = 144.25
$ws.Range("K2").Value = 104.2
$ws.Range("L2").Value = 144.25
$ws.Range("M2").Value = 8.799999999999997
$ws.Range("N2").Value = -370.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 31999
$ws.Range("J5").Value = 31999
$ws.Range("L5").Value = 31999
$ws.Range("N5").Value = -32225
$ws.Range("H16").Value = 1598.1482
$ws.Range("I16").Value = 1418.2916
$ws.Range("K16").Value = 1418.2916
$ws.Range("M16").Value = -1248.2916
$ws.Range("H20").Value = 122300000
$ws.Range("J20").Value = 550000000
$ws.Range("L20").Value = 550000000
$ws.Range("N20").Value = -550000452
$ws.Range("H22").Value = 5620.92
$ws.Range("I22").Value = 2958.7
$ws.Range("K22").Value = 2958.7
$ws.Range("M22").Value = -2663.7
$ws.Range("H27").Value = 5620.92
$ws.Range("I27").Value = 2958.7
$ws.Range("K27").Value = 2958.7
$ws.Range("M27").Value = -2851.7
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("H68").Value = 4599.6665
$ws.Range("I68").Value = 4599.6665
$ws.Range("K68").Value = 4599.6665
$ws.Range("M68").Value = -3850.6665
$ws.Range("H71").Value = 4599.6665
$ws.Range("I71").Value = 4599.6665
$ws.Range("K71").Value = 22998.3325
$ws.Range("M71").Value = -19254.3325
$ws.Range("H136").Value = 3519.7114
$ws.Range("J136").Value = 3718.2354
$ws.Range("L136").Value = 11154.7062
$ws.Range("N136").Value = -16254.7062
$ws.Range("M29").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 40000000
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("H24").Value = 9500
$ws.Range("I24").Value = 10000
$ws.Range("J24").Value = 9000
$ws.Range("K24").Value = 10000
$ws.Range("L24").Value = 9000
$ws.Range("M24").Value = -9770
$ws.Range("N24").Value = -9460
$ws.Range("H28").Value = 10000
$ws.Range("J28").Value = 10000
$ws.Range("L28").Value = 10000
$ws.Range("N28").Value = -10696
$ws.Range("H35").Value = 40000000
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("H122").Value = 1084.3158
$ws.Range("I122").Value = 980.34375
$ws.Range("J122").Value = 1638.8334
$ws.Range("K122").Value = 2941.03125
$ws.Range("L122").Value = 4916.5002
$ws.Range("M122").Value = -491.03125
$ws.Range("N122").Value = -9816.5002
$ws.Range("H124").Value = 82000
$ws.Range("J124").Value = 82000
$ws.Range("L124").Value = 82000
$ws.Range("N124").Value = -91820
$ws.Range("H132").Value = 1891.8246
$ws.Range("I132").Value = 1042.6765
$ws.Range("J132").Value = 3147.087
$ws.Range("K132").Value = 3128.0295
$ws.Range("L132").Value = 9441.261
$ws.Range("M132").Value = -598.0295000000001
$ws.Range("N132").Value = -14501.261
$ws.Range("N21").ClearContents()
$ws.Range("N35").ClearContents()
